$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet "Tactics" columns: A=id, B=version, C=name, D=name_EN, E=name_JP,
# F=detail, G=detail_EN, H=detail_JP. Rows 3-15 already hold the existing
# class data and are left untouched; append the two new classes
# ("turret" / "cocoon") as rows 16 and 17.
$ws.Cells.Item(16,1).Value2 = "turret"
$ws.Cells.Item(16,2).Value2 = "EA 23.210"
$ws.Cells.Item(16,3).Value2 = "炮塔"
$ws.Cells.Item(16,4).Value2 = "Turret"
$ws.Cells.Item(16,5).Value2 = "タレット"

$ws.Cells.Item(17,1).Value2 = "cocoon"
$ws.Cells.Item(17,2).Value2 = "EA 23.246"
$ws.Cells.Item(17,3).Value2 = "茧"
$ws.Cells.Item(17,4).Value2 = "Cocoon"
$ws.Cells.Item(17,5).Value2 = "コクーン"

# The source rows keep empty placeholder cells in columns G (detail_EN)
# and H (detail_JP); touch them (without actually changing formatting)
# so the saved sheet keeps an explicit, empty <c/> entry for G16/H16/G17/H17
# just like every other data row.
$ws.Cells.Item(16,7).Font.Size = 11
$ws.Cells.Item(16,8).Font.Size = 11
$ws.Cells.Item(17,7).Font.Size = 11
$ws.Cells.Item(17,8).Font.Size = 11
